# GitHub Actions refresh of the cryptos list (Thu Aug 17 03:33:47 UTC 2023).
# Column D ("Price") holds every figure as literal Text in the source sheet -
# even rows like "28.714.95" that are not parseable as a single Excel number -
# so each new price is entered with a leading apostrophe to force text entry
# (otherwise Excel would silently coerce it to a Number, e.g. dropping the
# trailing zero in "1.420"), then the quote-prefix style flag that apostrophe
# entry leaves behind is cleared by resetting the cell back to the Normal style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.713.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.63%  '
$ws.Range("D3").Value = '''1.803.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").Value = '''231.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("D6").Value = '''0.5947'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.06%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").Value = '''0.2782'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.05%  '
$ws.Range("D9").Value = '''0.06859'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.21%  '
$ws.Range("D10").Value = '''23.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("D11").Value = '''0.07549'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("D12").Value = '''1.807.12'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("D13").Value = '''4.729'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.84%  '
$ws.Range("D14").Value = '''0.6278'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").Value = '''2.049.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").Value = '''0.000009334'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.13%  '
$ws.Range("D17").Value = '''75.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.64%  '
$ws.Range("D18").Value = '''28.704.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.69%  '
$ws.Range("D19").Value = '''5.485'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.92%  '
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").Value = '''211.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.96%  '
$ws.Range("D22").Value = '''11.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.69%  '
$ws.Range("E23").Value = '  -2.05%  '
$ws.Range("D24").Value = '''1.004'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.32%  '
$ws.Range("D25").Value = '''154.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("D26").Value = '''7.858'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("D27").Value = '''0.1277'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("E28").Value = '  -0.68%  '
$ws.Range("D29").Value = '''1.446'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.09%  '
$ws.Range("D30").Value = '''0.06215'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.33%  '
$ws.Range("D31").Value = '''1.420'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.42%  '
$ws.Range("D32").Value = '''3.784'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("D33").Value = '''3.769'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.96%  '
$ws.Range("D34").Value = '''1.719'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("E35").Value = '  -6.11%  '
$ws.Range("D36").Value = '''0.6426'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.11%  '
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("D38").Value = '''2.725'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = '''0.01708'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("D40").Value = '''6.429'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("D41").Value = '''1.141.83'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.92%  '
$ws.Range("D42").Value = '''0.8668'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.08%  '
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("D44").Value = '''100.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").Value = '''1.968.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").Value = '''60.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.71%  '
$ws.Range("E47").Value = '  -5.58%  '
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").Value = '''8.358'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.11%  '
$ws.Range("D50").Value = '''0.05467'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.09%  '
$ws.Range("D51").Value = '''0.4493'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.49%  '
